# Scheduled Sheets update: refresh crafting-profit figures (H/I/J/K/L/M/N)
# across ALC/ARM/CRP/CUL/GSM/LTW/WVR per latest market-board pull.
# Cells with no remaining value are cleared entirely (matches source rows
# where a column legitimately has no data point this run).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1212.7667
$ws.Range("J17").Value = 1365.9584
$ws.Range("L17").Value = 4097.8752
$ws.Range("N17").Value = -4433.8752
$ws.Range("H132").Value = 2294.2
$ws.Range("I132").Value = 1061.5
$ws.Range("K132").Value = 3184.5
$ws.Range("M132").Value = -654.5
$ws.Range("H138").Value = 14930.741
$ws.Range("J138").Value = 14927.56
$ws.Range("L138").Value = 44782.68
$ws.Range("N138").Value = -55062.68

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7247.067
$ws.Range("I32").Value = 7550.4287
$ws.Range("K32").Value = 7550.4287
$ws.Range("M32").Value = -7263.4287
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
$ws.Range("H132").Value = 2753.8572
$ws.Range("I132").Value = 2202.6428
$ws.Range("K132").Value = 6607.928400000001
$ws.Range("M132").Value = -4077.928400000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3999.5
$ws.Range("I31").Value = 4500
$ws.Range("J31").Value = 3499
$ws.Range("K31").Value = 4500
$ws.Range("L31").Value = 3499
$ws.Range("M31").Value = -4205
$ws.Range("N31").Value = -4089
$ws.Range("H34").Value = 3999.5
$ws.Range("I34").Value = 4500
$ws.Range("J34").Value = 3499
$ws.Range("K34").Value = 4500
$ws.Range("L34").Value = 3499
$ws.Range("M34").Value = -4298
$ws.Range("N34").Value = -3903
$ws.Range("H58").Value = 3794.4
$ws.Range("I58").Value = 4022.5
$ws.Range("K58").Value = 4022.5
$ws.Range("M58").Value = -3819.5
$ws.Range("H132").Value = 4002.4
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 4003
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 12009
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -17069
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -12465
$ws.Range("N134").Value = ""
$ws.Range("H136").Value = 3794.4
$ws.Range("I136").Value = 4022.5
$ws.Range("K136").Value = 12067.5
$ws.Range("M136").Value = -9517.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 560.8
$ws.Range("I34").Value = 560.8
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1682.4
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1598.4
$ws.Range("N34").Value = ""
$ws.Range("H39").Value = 2000
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588
$ws.Range("H55").Value = 2400
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = ""
$ws.Range("H68").Value = 1062.25
$ws.Range("I68").Value = 775
$ws.Range("J68").Value = 1349.5
$ws.Range("K68").Value = 2325
$ws.Range("L68").Value = 4048.5
$ws.Range("M68").Value = -1514
$ws.Range("N68").Value = -5670.5
$ws.Range("H71").Value = 1062.25
$ws.Range("I71").Value = 775
$ws.Range("J71").Value = 1349.5
$ws.Range("K71").Value = 6975
$ws.Range("L71").Value = 12145.5
$ws.Range("M71").Value = -2919
$ws.Range("N71").Value = -20257.5
$ws.Range("H122").Value = 1449.8889
$ws.Range("J122").Value = 1364.6666
$ws.Range("L122").Value = 12281.9994
$ws.Range("N122").Value = -17181.9994

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12633.333
$ws.Range("I5").Value = 7950
$ws.Range("J5").Value = 22000
$ws.Range("K5").Value = 7950
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = -7838
$ws.Range("N5").Value = -22224
$ws.Range("H45").Value = 44000
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 78000
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 78000
$ws.Range("M45").Value = -9441
$ws.Range("N45").Value = -79118
$ws.Range("H113").Value = 1962.909
$ws.Range("I113").Value = 1919.4
$ws.Range("J113").Value = 1999.1666
$ws.Range("K113").Value = 1919.4
$ws.Range("L113").Value = 1999.1666
$ws.Range("M113").Value = 250.5999999999999
$ws.Range("N113").Value = -6339.1666

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2722.25
$ws.Range("I7").Value = 2598
$ws.Range("K7").Value = 2598
$ws.Range("M7").Value = -2486
$ws.Range("H106").Value = 55465.332
$ws.Range("J106").Value = 55465.332
$ws.Range("L106").Value = 55465.332
$ws.Range("N106").Value = -57989.332
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 2722.25
$ws.Range("I126").Value = 2598
$ws.Range("K126").Value = 7794
$ws.Range("M126").Value = -5324
$ws.Range("H136").Value = 3471.4285
$ws.Range("I136").Value = 3450
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 10350
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -7800
$ws.Range("N136").Value = -15600

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 59999.75
$ws.Range("J119").Value = 59999.75
$ws.Range("L119").Value = 59999.75
$ws.Range("N119").Value = -69675.75
$ws.Range("H122").Value = 5333
$ws.Range("I122").Value = 5333
$ws.Range("K122").Value = 15999
$ws.Range("M122").Value = -13549
$ws.Range("H132").Value = 3313.4666
$ws.Range("I132").Value = 2967
$ws.Range("K132").Value = 8901
$ws.Range("M132").Value = -6371
$ws.Range("H136").Value = 3642.8
$ws.Range("I136").Value = 3553.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 10660.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -8110.5
$ws.Range("N136").Value = -17100

